$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.243.13'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '1.800.22'
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'339.00"
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = "'0.9989"
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').Value = "'0.4780"
$ws.Range('E7').Value = '  +27.22%  '
$ws.Range('D8').Value = "'0.3629"
$ws.Range('E8').Value = '  +8.45%  '
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.07682"
$ws.Range('E10').Value = '  +7.57%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = "'1.144"
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('D12').Value = "'22.52"
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').Value = "'1.000"
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = "'6.272"
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').Value = "'7.283"
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '1.797.11'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = "'0.00001091"
$ws.Range('E17').Value = '  +4.08%  '
$ws.Range('D18').Value = "'0.06712"
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').Value = "'81.64"
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('D21').Value = "'17.29"
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('D22').Value = "'6.404"
$ws.Range('E22').Value = '  +2.37%  '
$ws.Range('D23').Value = '28.234.06'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('D24').Value = "'11.99"
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').Value = "'2.403"
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').Value = "'20.55"
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D27').Value = "'2.405"
$ws.Range('E27').Value = '  +4.14%  '
$ws.Range('D28').Value = "'151.83"
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = '2.002.33'
$ws.Range('E29').Value = '  +2.32%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'133.64"
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'1.274"
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('D32').Value = "'4.066"
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').Value = "'5.925"
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('D34').Value = "'0.09580"
$ws.Range('E34').Value = '  +10.73%  '
$ws.Range('D35').Value = "'0.02374"
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').Value = "'12.14"
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = "'0.06298"
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('D38').Value = "'0.6650"
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('D39').Value = "'5.205"
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').Value = "'0.2170"
$ws.Range('E40').Value = '  +3.20%  '
$ws.Range('E41').Value = '  +2.15%  '
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('D43').Value = "'8.091"
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = "'0.9986"
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  +3.15%  '
$ws.Range('D46').Value = "'3.870"
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').Value = "'0.6115"
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = "'128.35"
$ws.Range('E48').Value = '  -0.44%  '
$ws.Range('D49').Value = "'2.034"
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').Value = "'1.173"
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = "'0.07099"
$ws.Range('E51').Value = '  -1.08%  '
